# Auto-generated edit script applying the cryptos.xlsx diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '63.669.70'
$ws.Range('E2').Value = '  -3.05%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.145.01'
$ws.Range('E3').Value = '  -3.84%  '
$ws.Range('E4').Value = '  -0.06%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '563.93'
$ws.Range('E5').Value = '  -3.16%  '
$ws.Range('E6').Value = '  -6.89%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.599'
$ws.Range('E7').Value = '  -6.57%  '
$ws.Range('E8').Value = '  +0.09%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '3.141.60'
$ws.Range('E9').Value = '  -3.95%  '
$ws.Range('E10').Value = '  -4.42%  '
$ws.Range('E11').Value = '  -1.33%  '
$ws.Range('E12').Value = '  -5.22%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '3.690.23'
$ws.Range('E13').Value = '  -3.90%  '
$ws.Range('E14').Value = '  -1.75%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '63.775.78'
$ws.Range('E15').Value = '  -3.12%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '25.00'
$ws.Range('E16').Value = '  -3.51%  '
$ws.Range('E17').Value = '  -3.23%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '3.143.94'
$ws.Range('E18').Value = '  -2.80%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '409.12'
$ws.Range('E19').Value = '  -3.79%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '12.65'
$ws.Range('E20').Value = '  -3.61%  '
$ws.Range('E21').Value = '  -3.50%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '7.02'
$ws.Range('E22').Value = '  -4.28%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.999'
$ws.Range('E23').Value = '  -0.15%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '69.62'
$ws.Range('E24').Value = '  -2.92%  '
$ws.Range('E25').Value = '  +1.74%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.485'
$ws.Range('E26').Value = '  -4.20%  '
$ws.Range('E27').Value = '  -6.58%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '8.60'
$ws.Range('E28').Value = '  -2.42%  '
$ws.Range('E29').Value = '  -0.06%  '
$ws.Range('B30').Value = 'USDe'
$ws.Range('C30').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.999'
$ws.Range('E30').Value = '  -0.08%  '
$ws.Range('B31').Value = 'PancakeSwap'
$ws.Range('C31').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.81'
$ws.Range('E31').Value = '  -6.89%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '21.52'
$ws.Range('E32').Value = '  -2.86%  '
$ws.Range('E33').Value = '  -3.66%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '6.25'
$ws.Range('E34').Value = '  -4.62%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.11'
$ws.Range('E35').Value = '  -5.56%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '154.49'
$ws.Range('E36').Value = '  -2.99%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.35'
$ws.Range('E37').Value = '  -4.42%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.684.30'
$ws.Range('E38').Value = '  -3.51%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.67'
$ws.Range('E39').Value = '  -6.50%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '24.16'
$ws.Range('E40').Value = '  -7.94%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '4.11'
$ws.Range('E41').Value = '  -4.49%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '38.54'
$ws.Range('E42').Value = '  -3.48%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.701'
$ws.Range('E43').Value = '  -8.18%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0616'
$ws.Range('E44').Value = '  -5.98%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '5.54'
$ws.Range('E45').Value = '  -5.16%  '
$ws.Range('B46').Value = 'Bittensor'
$ws.Range('C46').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '292.18'
$ws.Range('E46').Value = '  -6.86%  '
$ws.Range('B47').Value = 'VeChain'
$ws.Range('C47').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0258'
$ws.Range('E47').Value = '  -3.02%  '
$ws.Range('B48').Value = 'InjectiveProtocol'
$ws.Range('C48').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '21.49'
$ws.Range('E48').Value = '  -6.37%  '
$ws.Range('E49').Value = '  -11.69%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.999'
$ws.Range('E50').Value = '  -0.16%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0983'
$ws.Range('E51').Value = '  -5.57%  '
